$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 45.8
$ws.Range("D2").Value = 14.61
$ws.Range("E2").Value = 0.42
$ws.Range("F2").Value = 20.76
$ws.Range("G2").Value = 0.04
$ws.Range("H2").Value = 0.01
$ws.Range("I2").Value = 0.07000000000000001
$ws.Range("J2").Value = 0.07000000000000001
$ws.Range("L2").Value = 0.1
$ws.Range("N2").Value = 0.03
$ws.Range("O2").Value = 0.64
$ws.Range("P2").Value = 0.16
$ws.Range("Q2").Value = 0.04
$ws.Range("T2").Value = 0.26
$ws.Range("Z2").Value = 0.28
$ws.Range("AA2").Value = 0.1
$ws.Range("C3").Value = 68.54000000000001
$ws.Range("D3").Value = 47.91
$ws.Range("F3").Value = 18.41
$ws.Range("I3").Value = 0.39
$ws.Range("K3").Value = 0.13
$ws.Range("M3").Value = 0.13
$ws.Range("O3").Value = 0.39
$ws.Range("S3").Value = 0.13
$ws.Range("T3").Value = 0.39
$ws.Range("W3").Value = 4.57
$ws.Range("C4").Value = 59.88
$ws.Range("D4").Value = 32.56
$ws.Range("F4").Value = 50.43
$ws.Range("G4").Value = 0.25
$ws.Range("J4").Value = 0.16
$ws.Range("K4").Value = 0.63
$ws.Range("O4").Value = 0.7
$ws.Range("S4").Value = 0.03
$ws.Range("T4").Value = 0.16
$ws.Range("Z4").Value = 0.19
$ws.Range("C5").Value = 48.7
$ws.Range("D5").Value = 23.25
$ws.Range("E5").Value = 0.09
$ws.Range("F5").Value = 10.49
$ws.Range("H5").Value = 0.03
$ws.Range("I5").Value = 0.04
$ws.Range("L5").Value = 0.05
$ws.Range("M5").Value = 0.45
$ws.Range("N5").Value = 0.03
$ws.Range("O5").Value = 0.39
$ws.Range("P5").Value = 0.05
$ws.Range("T5").Value = 0.25
$ws.Range("Z5").Value = 0.19
$ws.Range("C6").Value = 52.1
$ws.Range("D6").Value = 11.77
$ws.Range("F6").Value = 22.23
$ws.Range("G6").Value = 0.03
$ws.Range("H6").Value = 0.05
$ws.Range("I6").Value = 0.15
$ws.Range("J6").Value = 0.01
$ws.Range("K6").Value = 0.6
$ws.Range("L6").Value = 0.05
$ws.Range("M6").Value = 0.92
$ws.Range("N6").Value = 0.07000000000000001
$ws.Range("O6").Value = 0.75
$ws.Range("P6").Value = 0.33
$ws.Range("Q6").Value = 0.01
$ws.Range("S6").Value = 0.02
$ws.Range("T6").Value = 0.14
$ws.Range("U6").Value = 0.03
$ws.Range("Y6").Value = 0.02
$ws.Range("Z6").Value = 0.41
$ws.Range("AA6").Value = 0.05
$ws.Range("C7").Value = 68.37
$ws.Range("D7").Value = 58.02
$ws.Range("F7").Value = 32.51
$ws.Range("I7").Value = 0.15
$ws.Range("J7").Value = 0.15
$ws.Range("K7").Value = 1.02
$ws.Range("M7").Value = 0.87
$ws.Range("O7").Value = 0.44
$ws.Range("P7").Value = 0.15
$ws.Range("T7").Value = 0.44
$ws.Range("Z7").Value = 0.29
$ws.Range("AA7").Value = 0.29
$ws.Range("C8").Value = 63.61
$ws.Range("D8").Value = 29.75
$ws.Range("E8").Value = 0.79
$ws.Range("F8").Value = 29.51
$ws.Range("I8").Value = 0.08
$ws.Range("K8").Value = 0.4
$ws.Range("L8").Value = 0.16
$ws.Range("M8").Value = 0.47
$ws.Range("O8").Value = 0.4
$ws.Range("Z8").Value = 0.08
$ws.Range("C9").Value = 42.73
$ws.Range("D9").Value = 23.12
$ws.Range("F9").Value = 18.48
$ws.Range("G9").Value = 0.02
$ws.Range("H9").Value = 0.02
$ws.Range("I9").Value = 0.08
$ws.Range("J9").Value = 0.04
$ws.Range("L9").Value = 0.03
$ws.Range("O9").Value = 0.45
$ws.Range("P9").Value = 0.11
$ws.Range("S9").Value = 0
$ws.Range("Y9").Value = 0.03
$ws.Range("Z9").Value = 0.13
$ws.Range("AA9").Value = 0.04
$ws.Range("C10").Value = 56.27
$ws.Range("D10").Value = 44.09
$ws.Range("E10").Value = 0.05
$ws.Range("F10").Value = 27.63
$ws.Range("G10").Value = 0.04
$ws.Range("H10").Value = 0.03
$ws.Range("I10").Value = 0.02
$ws.Range("J10").Value = 0.01
$ws.Range("L10").Value = 0.02
$ws.Range("O10").Value = 0.6899999999999999
$ws.Range("P10").Value = 0.06
$ws.Range("S10").Value = 0.01
$ws.Range("Y10").Value = 0.01
$ws.Range("Z10").Value = 0.14
$ws.Range("AA10").Value = 0.15
$ws.Range("C11").Value = 50.9
$ws.Range("D11").Value = 43.85
$ws.Range("F11").Value = 24.82
$ws.Range("G11").Value = 0.01
$ws.Range("H11").Value = 0.01
$ws.Range("I11").Value = 0.01
$ws.Range("J11").Value = 0.01
$ws.Range("L11").Value = 0.11
$ws.Range("M11").Value = 0.73
$ws.Range("N11").Value = 0.02
$ws.Range("O11").Value = 0.59
$ws.Range("P11").Value = 0.28
$ws.Range("Q11").Value = 0.03
$ws.Range("S11").Value = 0.03
$ws.Range("U11").Value = 0.03
$ws.Range("Y11").Value = 0.01
$ws.Range("Z11").Value = 0.31
$ws.Range("AA11").Value = 0.04
$ws.Range("C12").Value = 65.06999999999999
$ws.Range("D12").Value = 56.45
$ws.Range("E12").Value = 0.27
$ws.Range("F12").Value = 37.91
$ws.Range("G12").Value = 0.06
$ws.Range("I12").Value = 0.15
$ws.Range("J12").Value = 0.06
$ws.Range("K12").Value = 0.31
$ws.Range("L12").Value = 0.1
$ws.Range("O12").Value = 0.52
$ws.Range("P12").Value = 0.19
$ws.Range("Q12").Value = 0.08
$ws.Range("S12").Value = 0.08
$ws.Range("T12").Value = 0.19
$ws.Range("Z12").Value = 0.57
$ws.Range("AA12").Value = 0.02
$ws.Range("C13").Value = 54.94
$ws.Range("D13").Value = 9.4
$ws.Range("E13").Value = 0.11
$ws.Range("F13").Value = 19.87
$ws.Range("I13").Value = 0.1
$ws.Range("J13").Value = 0.03
$ws.Range("L13").Value = 0.08
$ws.Range("M13").Value = 0.33
$ws.Range("O13").Value = 0.48
$ws.Range("P13").Value = 0.04
$ws.Range("S13").Value = 0.01
$ws.Range("Y13").Value = 0.04
$ws.Range("Z13").Value = 0.07000000000000001
$ws.Range("C14").Value = 53.91
$ws.Range("D14").Value = 28.75
$ws.Range("F14").Value = 31.58
$ws.Range("G14").Value = 0.03
$ws.Range("I14").Value = 0.09
$ws.Range("J14").Value = 0.01
$ws.Range("M14").Value = 0.47
$ws.Range("O14").Value = 0.79
$ws.Range("P14").Value = 0.09
$ws.Range("Q14").Value = 0.01
$ws.Range("Z14").Value = 0.14
$ws.Range("AA14").Value = 0.07000000000000001
$ws.Range("C15").Value = 37.05
$ws.Range("D15").Value = 45.01
$ws.Range("F15").Value = 9.06
$ws.Range("I15").Value = 0.03
$ws.Range("J15").Value = 0.03
$ws.Range("K15").Value = 0.43
$ws.Range("L15").Value = 0.01
$ws.Range("N15").Value = 0.01
$ws.Range("O15").Value = 0.47
$ws.Range("P15").Value = 0.05
$ws.Range("Y15").Value = 0.03
$ws.Range("Z15").Value = 0.19
$ws.Range("C16").Value = 52.17
$ws.Range("D16").Value = 66.95
$ws.Range("F16").Value = 37.62
$ws.Range("G16").Value = 0.04
$ws.Range("I16").Value = 0.15
$ws.Range("K16").Value = 0.92
$ws.Range("L16").Value = 0.04
$ws.Range("M16").Value = 1.62
$ws.Range("O16").Value = 0.59
$ws.Range("P16").Value = 0.41
$ws.Range("Q16").Value = 0.18
$ws.Range("S16").Value = 0.04
$ws.Range("T16").Value = 0.22
$ws.Range("Y16").Value = 0.04
$ws.Range("Z16").Value = 0.8100000000000001
$ws.Range("AA16").Value = 0.11
$ws.Range("C17").Value = 43.48
$ws.Range("D17").Value = 42.57
$ws.Range("F17").Value = 55.34
$ws.Range("H17").Value = 0.09
$ws.Range("K17").Value = 0.45
$ws.Range("M17").Value = 0.36
$ws.Range("O17").Value = 0.63
$ws.Range("P17").Value = 0.27
$ws.Range("T17").Value = 0.09
$ws.Range("Z17").Value = 0.18
